# Horarios Línea 141 - actualización de datos (scrape 07:50:28)
# Actualiza las 3 hojas (LP1912, LP1912-215, 6203-6173) con las nuevas filas
# scrapeadas, las horas "Última actualización" / "Total filas" y las filas
# existentes que cambiaron de posición/orden.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Hoja "LP1912"
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item('LP1912')
$rows = @{}
$rows[2] = @('Última actualización: 07:50:28', $null, $null, $null, $null)
$rows[3] = @('Total filas: 96', $null, $null, $null, $null)
$rows[43] = @('06:54:14', '07:11', '23_HERNANDEZ', 17, 'LP1912')
$rows[44] = @('06:54:14', '07:11', '215A_EL PATO', 17, 'LP1912')
$rows[49] = @('07:19:11', '07:20', '10_OLMOS', 1, 'LP1912')
$rows[50] = @('07:19:11', '07:20', '16_SANTA ANA', 1, 'LP1912')
$rows[53] = @('05:55:25', '07:31', '11_ETCHEVERRY', 96, 'LP1912')
$rows[54] = @('05:55:25', '07:31', '16_SANTA ANA', 96, 'LP1912')
$rows[64] = @('07:50:28', '07:50', '16_SANTA ANA', 0, 'LP1912')
$rows[65] = @('07:50:28', '07:51', '215D_EL PATO', 1, 'LP1912')
$rows[66] = @('07:19:11', '07:52', '215D_EL PATO', 33, 'LP1912')
$rows[67] = @('07:50:28', '07:55', '10_OLMOS', 5, 'LP1912')
$rows[68] = @('07:50:28', '07:58', '16_SANTA ANA', 8, 'LP1912')
$rows[69] = @('07:19:11', '08:00', '23_HERNANDEZ', 41, 'LP1912')
$rows[70] = @('06:26:08', '08:01', '23_HERNANDEZ', 95, 'LP1912')
$rows[71] = @('07:50:28', '08:02', '23_HERNANDEZ', 12, 'LP1912')
$rows[72] = @('07:50:28', '08:03', '11_ETCHEVERRY', 13, 'LP1912')
$rows[73] = @('07:19:11', '08:04', '11_ETCHEVERRY', 45, 'LP1912')
$rows[74] = @('06:54:14', '08:06', '23_HERNANDEZ', 72, 'LP1912')
$rows[75] = @('07:50:28', '08:12', '15_ABASTO', 22, 'LP1912')
$rows[76] = @('07:50:28', '08:13', '10_OLMOS', 23, 'LP1912')
$rows[77] = @('07:50:28', '08:21', '26_HERNANDEZ', 31, 'LP1912')
$rows[78] = @('07:50:28', '08:22', '16_P MOR-SANTA ANA', 32, 'LP1912')
$rows[79] = @('07:19:11', '08:23', '215B_EL PATO', 64, 'LP1912')
$rows[80] = @('07:50:28', '08:23', '215B_EL PATO', 33, 'LP1912')
$rows[81] = @('07:50:28', '08:27', '84_COLONIA URQUIZA-ESC 49', 37, 'LP1912')
$rows[82] = @('07:50:28', '08:37', '23_HERNANDEZ', 47, 'LP1912')
$rows[83] = @('07:19:11', '08:42', '81_EL PELIGRO', 83, 'LP1912')
$rows[84] = @('07:50:28', '08:43', '14_ABASTO', 53, 'LP1912')
$rows[85] = @('07:19:11', '08:44', '14_ABASTO', 85, 'LP1912')
$rows[86] = @('07:50:28', '08:54', '17_ROMERO', 64, 'LP1912')
$rows[87] = @('07:50:28', '09:01', '215A_EL PATO', 71, 'LP1912')
$rows[88] = @('07:19:11', '09:02', '215A_EL PATO', 103, 'LP1912')
$rows[89] = @('07:50:28', '09:03', '11_ETCHEVERRY', 73, 'LP1912')
$rows[90] = @('07:50:28', '09:10', '16_P MOR-SANTA ANA', 80, 'LP1912')
$rows[91] = @('07:19:11', '09:11', '16_P MOR-SANTA ANA', 112, 'LP1912')
$rows[92] = @('07:50:28', '09:11', '81_EL PELIGRO', 81, 'LP1912')
$rows[93] = @('07:50:28', '09:16', '27_EL RETIRO', 86, 'LP1912')
$rows[94] = @('07:19:11', '09:17', '27_EL RETIRO', 118, 'LP1912')
$rows[95] = @('07:50:28', '09:21', '26_HERNANDEZ', 91, 'LP1912')
$rows[96] = @('07:50:28', '09:22', '17_ROMERO', 92, 'LP1912')
$rows[97] = @('07:50:28', '09:23', '11_ETCHEVERRY', 93, 'LP1912')
$rows[98] = @('07:50:28', '09:32', '15_ABASTO', 102, 'LP1912')
$rows[99] = @('07:50:28', '09:33', '10_OLMOS', 103, 'LP1912')
$rows[100] = @('07:50:28', '09:42', '215C_EL PATO', 112, 'LP1912')
$rows[101] = @('07:50:28', '09:43', '14_ABASTO', 113, 'LP1912')
foreach ($r in $rows.Keys) {
    $row = $rows[$r]
    for ($c = 1; $c -le 5; $c++) {
        $val = $row[$c-1]
        if ($val -ne $null) {
            $ws.Cells.Item($r, $c).Value = $val
        }
    }
}

# ------------------------------------------------------------------
# Hoja "LP1912-215"
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item('LP1912-215')
$rows = @{}
$rows[2] = @('Última actualización: 07:50:28', $null, $null, $null, $null)
$rows[3] = @('Total filas: 16', $null, $null, $null, $null)
$rows[16] = @('07:50:28', '07:51', '215D_EL PATO', 1, 'LP1912')
$rows[18] = @('07:50:28', '08:23', '215B_EL PATO', 33, 'LP1912')
$rows[19] = @('07:50:28', '09:01', '215A_EL PATO', 71, 'LP1912')
$rows[20] = @('07:19:11', '09:02', '215A_EL PATO', 103, 'LP1912')
$rows[21] = @('07:50:28', '09:42', '215C_EL PATO', 112, 'LP1912')
foreach ($r in $rows.Keys) {
    $row = $rows[$r]
    for ($c = 1; $c -le 5; $c++) {
        $val = $row[$c-1]
        if ($val -ne $null) {
            $ws.Cells.Item($r, $c).Value = $val
        }
    }
}

# ------------------------------------------------------------------
# Hoja "6203-6173"
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item('6203-6173')
$rows = @{}
$rows[2] = @('Última actualización: 07:50:28', $null, $null, $null, $null)
$rows[3] = @('Total filas: 13', $null, $null, $null, $null)
$rows[15] = @('07:50:28', '08:09', '215C_LA PLATA', 19, 'L6203')
$rows[16] = @('06:54:14', '08:31', '215A_LA PLATA', 97, 'L6173')
$rows[17] = @('07:50:28', '08:35', '215A_LA PLATA', 45, 'L6173')
$rows[18] = @('07:50:28', '09:09', '215D_LA PLATA', 79, 'L6203')
foreach ($r in $rows.Keys) {
    $row = $rows[$r]
    for ($c = 1; $c -le 5; $c++) {
        $val = $row[$c-1]
        if ($val -ne $null) {
            $ws.Cells.Item($r, $c).Value = $val
        }
    }
}
